$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (rows 2-7)
$data = @(
    @("1326573", "https://aiesec.org/opportunity/global-talent/1326573", "[Impact Porto Alegre] - SBD Engineering Intern", "Rio Claro, SP, Brasil", "No", "0 applicants", "6 - 18 Months", "Whirlpool Corporation"),
    @("1326162", "https://aiesec.org/opportunity/global-talent/1326162", "Tech Sales Development Representative( swedish Only)", "Bournemouth, Royaume-Uni", "No", "6 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("1326160", "https://aiesec.org/opportunity/global-talent/1326160", "Tech Sales Development Representative( dutch  Only)", "Bournemouth, Royaume-Uni", "No", "3 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("1326159", "https://aiesec.org/opportunity/global-talent/1326159", "Tech Sales Development Representative( Spanish Only)", "Bournemouth, Royaume-Uni", "No", "27 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("1326156", "https://aiesec.org/opportunity/global-talent/1326156", "Tech Sales Development Representative( French/ Swiss/ Belgian Only)", "Bournemouth, Royaume-Uni", "No", "16 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("1326152", "https://aiesec.org/opportunity/global-talent/1326152", "Tech Sales Development Representative( German / Austrian Only)", "Bournemouth, Royaume-Uni", "No", "2 applicants", "6 - 18 Months", "EIMS Ltd")
)

$rowIndex = 2
foreach ($rowData in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $value = $rowData[$col - 1]
        $cell = $ws.Cells.Item($rowIndex, $col)
        if ($col -eq 1) {
            # Column A holds numeric-looking IDs that must stay text. A bare
            # assignment gets auto-coerced to a number, so use the
            # apostrophe quote-prefix trick to force text, then reset the
            # cell style back to Normal so no quote-prefix style lingers.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
    $rowIndex++
}

# Update column widths. Excel's COM ColumnWidth setter round-trips through a
# pixel conversion that adds ~0.8333 (5/6) to the stored width value, so we
# subtract that offset here to land on the exact target widths seen in the
# OOXML (<col width="..."/>).
$widthOffset = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 55 - $widthOffset
$ws.Columns.Item(3).ColumnWidth = 70 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 27 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 16 - $widthOffset
$ws.Columns.Item(7).ColumnWidth = 16 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 24 - $widthOffset
